$d = $word.ActiveDocument

# 1. Remove the leftover "_GoBack" bookmark (left after "Beginning roll competition").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Merge the split runs into single runs by replacing the paragraph text
#    via Find/Replace (Word collapses the matched range into one run).
$find = $d.Content.Find
$find.Execute("Players earn cards when hexes are rolled", $true, $false, $false, $false, $false, $true, 1, $false, "Players earn cards when hexes are rolled", 2) | Out-Null

$find = $d.Content.Find
$find.Execute("Players use cards when building and then receive Victory Points", $true, $false, $false, $false, $false, $true, 1, $false, "Players use cards when building and then receive Victory Points", 2) | Out-Null

$find = $d.Content.Find
$find.Execute("The game ends when a player gets 10", $true, $false, $false, $false, $false, $true, 1, $false, "The game ends when a player gets 10", 2) | Out-Null

# 3. Fill in the lines-of-code numbers.
$find = $d.Content.Find
$find.Execute("Davis: ", $true, $false, $false, $false, $false, $true, 1, $false, "Davis: ~475", 2) | Out-Null

$find = $d.Content.Find
$find.Execute("Ethan: ", $true, $false, $false, $false, $false, $true, 1, $false, "Ethan: ~400", 2) | Out-Null

$find = $d.Content.Find
$find.Execute("Andrew: ", $true, $false, $false, $false, $false, $true, 1, $false, "Andrew: ~650", 2) | Out-Null
